$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price & 1h volume change columns)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.900.49"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.031.99"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.54"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.04"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.394"
$ws.Range("E9").Value = "  +6.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.12"
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  +6.57%  "
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.50"
$ws.Range("E13").Value = "  +18.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.892"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.38"
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.327.64"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.61"
$ws.Range("E17").Value = "  +3.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.027.22"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.789.89"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.62"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0886"
$ws.Range("E21").Value = "  +3.24%  "
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.16"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  -5.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.94"
$ws.Range("E27").Value = "  +4.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.140"
$ws.Range("E28").Value = "  +24.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "160.56"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.39"
$ws.Range("E30").Value = "  +3.94%  "
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0625"
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.53"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.57"
$ws.Range("E36").Value = "  +11.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.40"
$ws.Range("E37").Value = "  -4.68%  "
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.23"
$ws.Range("E40").Value = "  +25.44%  "
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("E42").Value = "  +3.15%  "
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "94.14"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.73"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.366.11"
$ws.Range("E49").Value = "  -4.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.91"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.219.34"
$ws.Range("E51").Value = "  +1.07%  "

# Row 32/33 swap: ImmutableX moves up to row 32 (was Filecoin), Filecoin moves to row 33
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.20"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.07"
$ws.Range("E33").Value = "  -1.22%  "

# Row 45/46 swap: VeChain moves up to row 45 (was InjectiveProtocol), InjectiveProtocol moves to row 46
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0216"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.03"
$ws.Range("E46").Value = "  +2.70%  "
